# "Group status attacks work" - duplicate the "Status Effect" sheet into a
# new "Group Status" sheet, retarget the former's stray status-icon cells
# onto the new sheet, and rename the "ParaNail" enemy entry to "Mushroom".

$wb = $excel.ActiveWorkbook
$statusEffect = $wb.Worksheets.Item("Status Effect")

# Duplicate "Status Effect" right after itself, becoming the new last sheet.
$statusEffect.Copy($null, $statusEffect)
$groupStatus = $wb.Worksheets.Item($statusEffect.Index + 1)
$groupStatus.Name = "Group Status"

# Row 3 on the original sheet no longer carries the stray L3/M3 status-icon
# reference cells.
$statusEffect.Range("L3").ClearContents()
$statusEffect.Range("M3").ClearContents()

# The copied sheet also loses those two cells, but gains a proper Enemy
# row: the "ParaNail" entry is being renamed to "Mushroom" (it is the only
# remaining reference to that shared string, so editing it in place renames
# it instead of allocating a new string) and the row gets its CLASS/LIVES
# values filled in like the other enemy rows.
$groupStatus.Range("L3").ClearContents()
$groupStatus.Range("M3").ClearContents()
$groupStatus.Range("B3").Value = "Mushroom"
$groupStatus.Range("C3").Value = "Enemy"
$groupStatus.Range("D3").Value = 1

# Restore the selections shown in each sheet's view, and make the new sheet
# the active tab (it was just created/activated by the copy).
[void]$statusEffect.Range("M3").Select()
[void]$groupStatus.Activate()
[void]$groupStatus.Range("B4").Select()
